# Update the "F" column (想去人数 / want-to-go count) values on the
# "展览" and "全部类型" worksheets, which contain the same data.

$wb = $excel.ActiveWorkbook

# Row -> new value for column F
$updates = @{
    2  = 138
    6  = 132
    7  = 1279
    8  = 1542
    10 = 403
    12 = 157
    14 = 66
    18 = 326
    19 = 1741
    22 = 176
    23 = 672
    26 = 4204
    27 = 13
    28 = 276
    29 = 1095
    30 = 490
    32 = 563
    34 = 264
    36 = 141
    37 = 14
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
